$wb = $excel.ActiveWorkbook

# Remove the "SortedData" worksheet entirely.
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("SortedData").Delete()

# Update the manual data entry rows on RawData.
$ws = $wb.Worksheets.Item("RawData")

$ws.Range("A2").Value = 111109
$ws.Range("B2").Value = "a"
$ws.Range("C2").Value = "b"
$ws.Range("D2").Value = "c"
$ws.Range("E2").Value = "N"
$ws.Range("F2").Value = "Y"
$ws.Range("G2").Value = "N"

$ws.Range("A3").Value = 111110
$ws.Range("B3").Value = "a"
$ws.Range("C3").Value = "b"
$ws.Range("D3").Value = "a"
$ws.Range("E3").Value = "Y"
$ws.Range("F3").Value = "N"
$ws.Range("G3").Value = "N"

[void]$ws.Range("G32").Select()
